$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G8").Value  = "Vijayapura (Bijapur)"
$ws.Range("G14").Value = "Vijayapura (Bijapur)"
$ws.Range("G21").Value = "Vijayapura (Bijapur)"
$ws.Range("G23").Value = "Uttara Kannada (Karwar)"
$ws.Range("G24").Value = "Vijayapura (Bijapur)"
$ws.Range("G25").Value = "Vijayapura (Bijapur)"
$ws.Range("G28").Value = "Vijayapura (Bijapur)"
$ws.Range("G33").Value = "Haveri"
